$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").Value = "'5.057"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '3HuobiTokenHT'
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").Value = "'0.05614"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '4CronosCRO'
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = "'6.539"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = "'3.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.8085"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = "'0.8430"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1338"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03250"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09413"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0005978"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14OneONE'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.006168"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15TigerCashTCH'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.501"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16LEOLEO'
$ws.Range("D18").Value = "'2.091"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.3183"
$ws.Range("D19").Style = "Normal"
$ws.Range("B20").Value = 'MandalaExchangeToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D20").Value = "'0.06876"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19MandalaExchangeTokenMDX'
$ws.Range("D22").Value = "'3.741"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.04676"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.1370"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.001243"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004531"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009698"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26NitroExNTXBestin24h'
$ws.Range("D40").Value = "'0.03642"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.1363"
$ws.Range("D41").Style = "Normal"
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.002723"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.003368"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("D44").Value = "'0.008066"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005279"
$ws.Range("D45").Style = "Normal"

Write-Output "edits applied"
